# Updated cryptos list on Fri Oct 27 02:24:39 UTC 2023 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures for each ranked coin, and
# keeps rows B-D39/40 and 45-47 in sync where the underlying ranking reordered
# a few coins (Aave <-> HuobiToken, Quant/Kaspa/FraxShare).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.997.32'
$ws.Range('E2').Value = '  -1.96%  '

$ws.Range('D3').Value = '1.786.53'
$ws.Range('E3').Value = '  -0.24%  '

$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').Value = "'221.21"
$ws.Range('E5').Value = '  -1.06%  '

$ws.Range('D6').Value = "'0.552"
$ws.Range('E6').Value = '  -0.10%  '

$ws.Range('D7').Value = "'0.998"
$ws.Range('E7').Value = '  -0.12%  '

$ws.Range('D8').Value = "'32.33"
$ws.Range('E8').Value = '  +0.44%  '

$ws.Range('D9').Value = "'0.283"
$ws.Range('E9').Value = '  +0.56%  '

$ws.Range('D10').Value = "'0.0711"
$ws.Range('E10').Value = '  +0.80%  '

$ws.Range('E11').Value = '  -0.72%  '

$ws.Range('D12').Value = '2.042.07'
$ws.Range('E12').Value = '  -0.29%  '

$ws.Range('D13').Value = '1.774.50'
$ws.Range('E13').Value = '  -1.33%  '

$ws.Range('D14').Value = "'10.79"
$ws.Range('E14').Value = '  -1.76%  '

$ws.Range('D15').Value = "'0.624"
$ws.Range('E15').Value = '  -1.26%  '

$ws.Range('D16').Value = '33.970.94'
$ws.Range('E16').Value = '  -1.95%  '

$ws.Range('D17').Value = "'4.15"
$ws.Range('E17').Value = '  -3.32%  '

$ws.Range('D18').Value = "'67.69"
$ws.Range('E18').Value = '  -1.97%  '

$ws.Range('D19').Value = "'243.50"
$ws.Range('E19').Value = '  -4.07%  '

$ws.Range('D20').Value = '0.0₃0781'
$ws.Range('E20').Value = '  -2.52%  '

$ws.Range('E21').Value = '  +0.06%  '

$ws.Range('D22').Value = "'10.81"
$ws.Range('E22').Value = '  +1.64%  '

$ws.Range('D23').Value = "'4.07"
$ws.Range('E23').Value = '  -2.68%  '

$ws.Range('D24').Value = "'2.11"
$ws.Range('E24').Value = '  -1.39%  '

$ws.Range('D25').Value = "'157.62"
$ws.Range('E25').Value = '  -1.94%  '

$ws.Range('D26').Value = "'16.32"
$ws.Range('E26').Value = '  -0.06%  '

$ws.Range('D27').Value = "'7.03"
$ws.Range('E27').Value = '  -1.00%  '

$ws.Range('E28').Value = '  -1.66%  '

$ws.Range('D29').Value = "'0.998"
$ws.Range('E29').Value = '  -0.03%  '

$ws.Range('D30').Value = "'0.0518"
$ws.Range('E30').Value = '  -1.06%  '

$ws.Range('E31').Value = '  +0.50%  '

$ws.Range('D32').Value = "'3.66"
$ws.Range('E32').Value = '  -2.94%  '

$ws.Range('D33').Value = "'3.49"
$ws.Range('E33').Value = '  -2.84%  '

$ws.Range('D34').Value = "'1.82"
$ws.Range('E34').Value = '  -2.50%  '

$ws.Range('D35').Value = '1.395.06'
$ws.Range('E35').Value = '  -2.85%  '

$ws.Range('D36').Value = "'0.639"
$ws.Range('E36').Value = '  +0.87%  '

$ws.Range('D37').Value = "'1.05"
$ws.Range('E37').Value = '  -0.46%  '

$ws.Range('D38').Value = "'0.0185"
$ws.Range('E38').Value = '  -3.36%  '

$ws.Range('B39').Value = 'Aave'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D39').Value = "'79.49"
$ws.Range('E39').Value = '  -6.07%  '

$ws.Range('B40').Value = 'HuobiToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D40').Value = "'2.35"
$ws.Range('E40').Value = '  +1.04%  '

$ws.Range('D41').Value = "'0.920"
$ws.Range('E41').Value = '  -0.27%  '

$ws.Range('D42').Value = "'2.71"
$ws.Range('E42').Value = '  -2.49%  '

$ws.Range('D43').Value = "'2.15"
$ws.Range('E43').Value = '  +1.88%  '

$ws.Range('E44').Value = '  -1.22%  '

$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = "'107.04"
$ws.Range('E45').Value = '  +1.74%  '

$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').Value = "'0.0492"
$ws.Range('E46').Value = '  +0.70%  '

$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = "'5.87"
$ws.Range('E47').Value = '  -1.42%  '

$ws.Range('D48').Value = '1.942.67'
$ws.Range('E48').Value = '  +0.16%  '

$ws.Range('D49').Value = "'0.998"
$ws.Range('E49').Value = '  -0.10%  '

$ws.Range('D50').Value = "'11.93"
$ws.Range('E50').Value = '  -0.30%  '

$ws.Range('E51').Value = '  +0.03%  '
